$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> FAPs)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5872273333333333
$ws.Range("H2").Value = 1.761682
$ws.Range("I2").Value = 0.07214749117712801
$ws.Range("J2").Value = 0.07214749117712801
$ws.Range("O2").Value = 0.4698491839234223
$ws.Range("P2").Value = 0.4698491839234222
$ws.Range("Q2").Value = 0.1335122022491111
$ws.Range("R2").Value = 1.201609820242
$ws.Range("S2").Value = 0.03389843985169591
$ws.Range("T2").Value = 0.0338984398516959

# Row 3 (ECs -> MuSCs)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5872273333333333
$ws.Range("H3").Value = 1.761682
$ws.Range("I3").Value = 0.07214749117712801
$ws.Range("J3").Value = 0.07214749117712801
$ws.Range("O3").Value = 0.3495359240395067
$ws.Range("P3").Value = 0.3495359240395067
$ws.Range("Q3").Value = 0.0993240226448889
$ws.Range("R3").Value = 0.8939162038040001
$ws.Range("S3").Value = 0.0252181399957296
$ws.Range("T3").Value = 0.0252181399957296

# Row 4 (ECs -> Resolving-Mac)
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5872273333333333
$ws.Range("H4").Value = 1.761682
$ws.Range("I4").Value = 0.07214749117712801
$ws.Range("J4").Value = 0.07214749117712801
$ws.Range("M4").Value = 0.08739966666666667
$ws.Range("N4").Value = 0.262199
$ws.Range("O4").Value = 0.180614892037071
$ws.Range("P4").Value = 0.180614892037071
$ws.Range("Q4").Value = 0.05132347319088889
$ws.Range("R4").Value = 0.461911258718
$ws.Range("S4").Value = 0.01303091132970251
$ws.Range("T4").Value = 0.01303091132970251

# Row 5 (FAPs -> FAPs)
$ws.Range("I5").Value = 0.927852508822872
$ws.Range("J5").Value = 0.927852508822872
$ws.Range("O5").Value = 0.4698491839234223
$ws.Range("P5").Value = 0.4698491839234222
$ws.Range("S5").Value = 0.4359507440717264
$ws.Range("T5").Value = 0.4359507440717263

# Row 6 (FAPs -> MuSCs)
$ws.Range("I6").Value = 0.927852508822872
$ws.Range("J6").Value = 0.927852508822872
$ws.Range("O6").Value = 0.3495359240395067
$ws.Range("P6").Value = 0.3495359240395067
$ws.Range("S6").Value = 0.3243177840437771
$ws.Range("T6").Value = 0.3243177840437771

# Row 7 (FAPs -> Resolving-Mac)
$ws.Range("I7").Value = 0.927852508822872
$ws.Range("J7").Value = 0.927852508822872
$ws.Range("M7").Value = 0.08739966666666667
$ws.Range("N7").Value = 0.262199
$ws.Range("O7").Value = 0.180614892037071
$ws.Range("P7").Value = 0.180614892037071
$ws.Range("Q7").Value = 0.6600453125217778
$ws.Range("R7").Value = 5.940407812696
$ws.Range("S7").Value = 0.1675839807073684
$ws.Range("T7").Value = 0.1675839807073684
